# Duplicate Leather data and test file
# Rolls the "Automation Item 65-69" / "Leather00-12" test rows forward to
# "Automation Item 70-74" / "Leather00-13" across the TestInventory,
# ItemCarousel and restocking sheets, and bumps the restocking date.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: TestInventory
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TestInventory")

$ws1.Range("D2").Value = "Leather00-13"
$ws1.Range("E2").Value = "LT-013"

$ws1.Range("D3").Value = "(Automation) Item 70"
$ws1.Range("E3").Value = "AT-IT-70"

$ws1.Range("D4").Value = "(Automation) Item 71"
$ws1.Range("E4").Value = "AT-IT-71"

$ws1.Range("D5").Value = "(Automation) Item 72"
$ws1.Range("E5").Value = "AT-IT-72"

$ws1.Range("D6").Value = "(Automation) Item 73"
$ws1.Range("E6").Value = "AT-IT-73"

$ws1.Range("D7").Value = "(Automation) Item 74"
$ws1.Range("E7").Value = "AT-IT-74"

# ---------------------------------------------------------------
# Sheet: ItemCarousel
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ItemCarousel")

$ws2.Range("C2").Value = "Leather00-13"
$ws2.Range("D2").Value = "LT-013"
$ws2.Range("E2").Value = "T0-LT-00-04"
$ws2.Range("F2").Value = "AT-LT-00-04"
$ws2.Range("G2").Value = "Auto Testing -004"

$ws2.Range("C3").Value = "(Automation) Item 70"
$ws2.Range("D3").Value = "AT-IT-70"
$ws2.Range("E3").Value = "T070"
$ws2.Range("F3").Value = "AT070"
$ws2.Range("G3").Value = "Auto Testing 70"

$ws2.Range("C4").Value = "(Automation) Item 71"
$ws2.Range("D4").Value = "AT-IT-71"
$ws2.Range("E4").Value = "T071"
$ws2.Range("F4").Value = "AT071"
$ws2.Range("G4").Value = "Auto Testing 71"

$ws2.Range("C5").Value = "(Automation) Item 72"
$ws2.Range("D5").Value = "AT-IT-72"
$ws2.Range("E5").Value = "T072"
$ws2.Range("F5").Value = "AT072"
$ws2.Range("G5").Value = "Auto Testing 72"

$ws2.Range("C6").Value = "(Automation) Item 73"
$ws2.Range("D6").Value = "AT-IT-73"
$ws2.Range("E6").Value = "T073"
$ws2.Range("F6").Value = "AT073"
$ws2.Range("G6").Value = "Auto Testing 73"

$ws2.Range("C7").Value = "(Automation) Item 74"
$ws2.Range("D7").Value = "AT-IT-74"
$ws2.Range("E7").Value = "T074"
$ws2.Range("F7").Value = "AT074"
$ws2.Range("G7").Value = "Auto Testing 74"

$ws2.Activate()
$ws2.Range("E13").Select()

# ---------------------------------------------------------------
# Sheet: restocking
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("restocking")

$ws3.Range("C2").Value = "Leather00-13"
$ws3.Range("D2").Value = "LT-013"
$ws3.Range("G2").Value = "2023-07-07"

$ws3.Range("C3").Value = "(Automation) Item 70"
$ws3.Range("D3").Value = "AT-IT-70"
$ws3.Range("G3").Value = "2023-07-07"

$ws3.Range("C4").Value = "(Automation) Item 71"
$ws3.Range("D4").Value = "AT-IT-71"
$ws3.Range("G4").Value = "2023-07-07"

$ws3.Range("C5").Value = "(Automation) Item 72"
$ws3.Range("D5").Value = "AT-IT-72"
$ws3.Range("G5").Value = "2023-07-07"

$ws3.Range("C6").Value = "(Automation) Item 73"
$ws3.Range("D6").Value = "AT-IT-73"
$ws3.Range("G6").Value = "2023-07-07"

$ws3.Range("C7").Value = "(Automation) Item 74"
$ws3.Range("D7").Value = "AT-IT-74"
$ws3.Range("G7").Value = "2023-07-07"

$ws3.Activate()
$ws3.Range("C2:D7").Select()

# Restore ItemCarousel as the active/selected sheet (matches the saved
# view state in the workbook: tabSelected stays on ItemCarousel).
$ws2.Activate()
